$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 222
$ws.Range("B3").Value = "castle"
$ws.Range("C3").Value = 256
$ws.Range("D3").Value = 896
$ws.Range("E3").Value = "MIME.png"

$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

$ws.Range("F9").Select() | Out-Null
